# Update the two-digit-divided-by-one-digit division problems in the
# first (and only) table of the worksheet. Each table cell containing a
# problem ("NN÷N=") is addressed directly by its (row, column) position so
# that cells sharing identical original text (e.g. the two "84÷6=" cells)
# can still be updated independently with distinct new values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "89÷8="  # was 19÷4=
$t.Cell(1, 2).Range.Text = "83÷6="  # was 67÷7=
$t.Cell(1, 3).Range.Text = "47÷8="  # was 82÷3=
$t.Cell(1, 4).Range.Text = "35÷8="  # was 70÷4=
$t.Cell(1, 5).Range.Text = "15÷3="  # was 99÷5=

# Row 5
$t.Cell(5, 1).Range.Text = "76÷7="  # was 71÷5=
$t.Cell(5, 2).Range.Text = "54÷9="  # was 84÷6=
$t.Cell(5, 3).Range.Text = "88÷3="  # was 58÷6=
$t.Cell(5, 4).Range.Text = "31÷5="  # was 33÷3=
$t.Cell(5, 5).Range.Text = "33÷8="  # was 84÷6=

# Row 9 (first cell "78÷6=" is unchanged)
$t.Cell(9, 2).Range.Text = "59÷2="  # was 61÷2=
$t.Cell(9, 3).Range.Text = "53÷9="  # was 92÷5=
$t.Cell(9, 4).Range.Text = "72÷6="  # was 51÷9=
$t.Cell(9, 5).Range.Text = "38÷4="  # was 31÷5=

# Row 13
$t.Cell(13, 1).Range.Text = "18÷3="  # was 59÷8=
$t.Cell(13, 2).Range.Text = "42÷8="  # was 20÷9=
$t.Cell(13, 3).Range.Text = "32÷9="  # was 99÷8=
$t.Cell(13, 4).Range.Text = "87÷5="  # was 52÷6=
$t.Cell(13, 5).Range.Text = "81÷9="  # was 35÷6=

# Row 17
$t.Cell(17, 1).Range.Text = "84÷7="  # was 30÷3=
$t.Cell(17, 2).Range.Text = "20÷2="  # was 52÷5=
$t.Cell(17, 3).Range.Text = "31÷6="  # was 48÷8=
$t.Cell(17, 4).Range.Text = "10÷5="  # was 25÷2=
$t.Cell(17, 5).Range.Text = "68÷6="  # was 87÷7=
